# Adding updated EMP 2020 data and new 2020 SMSCG data, and fixed EDI links
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("lat_long")
$ws.Activate()

# Insert a new row at row 358, pushing the existing row 358 ("606", the FMWT
# station previously right after "602") and everything below it down by one.
# This makes room for the new FMWT station "605".
$ws.Range("A358:D358").Insert()

# Give the new row the same look as the surrounding FMWT station rows:
# column A keeps the "FMWT" source-label style, columns B:D are plain
# numbers (no special formatting), like every other station row.
$ws.Cells.Item(357, 1).Copy()
$ws.Cells.Item(358, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B358:D358").Style = "Normal"

# New 2020 SMSCG station data.
$ws.Cells.Item(358, 1).Value = "FMWT"
$ws.Cells.Item(358, 2).Value = 605
$ws.Cells.Item(358, 3).Value = 38.148530000000001
$ws.Cells.Item(358, 4).Value = -122.05737999999999

# Leave the view where the editor ended up after making the change.
$ws.Application.ActiveWindow.ScrollRow = 334
$ws.Range("D358").Select()
